# Workbook: ya9-3-trabajo-infantil.xlsx
# Commit: "Datos que faltaban hasta el 10"
# Rename the sole sheet to "datos", add a new "metadatos" sheet describing
# the variables, their sources and the extraction date, and make the new
# sheet the active one (as left by the author).

$wb = $excel.ActiveWorkbook

# --- rename the original data sheet -----------------------------------
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

# --- add the metadata sheet right after it -----------------------------
$meta = $wb.Worksheets.Add($null, $datos)
$meta.Name = "metadatos"

# --- header row ----------------------------------------------------------
$meta.Range("A1").Value = "Variables"
$meta.Range("B1").Value = "Descripción"
$meta.Range("C1").Value = "Fuente"
$meta.Range("D1").Value = "Fecha_de_extracción"

# --- row 2: anno ---------------------------------------------------------
$meta.Range("A2").Value = "anno"
$meta.Range("B2").Value = "Año"
$meta.Range("C2").Value = "…"

# --- row 3: codmpio --------------------------------------------------------
$meta.Range("A3").Value = "codmpio"
$meta.Range("B3").Value = "Código del municipio"
$meta.Range("C3").Value = "…"

# --- row 4: SRPA_3 ---------------------------------------------------------
$meta.Range("A4").Value = "SRPA_3"
$meta.Range("B4").Value = "No. total de casos de niñas y adolescentes víctimas de trabajo infantil que ingresaron a proceso administrativo de restablecimiento de derechos`nNo. de adolescentes que ingresan al SRPA con una medida privativa de la libertad "
$meta.Range("C4").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- apply the shared font to the rows filled in so far (A1:D4) --------
$meta.Range("A1:D4").Font.Name = "Calibri"

# --- row 6: tasa (row 5 is filled afterwards, matching the source file) --
$meta.Range("A6").Value = "tasa"
$meta.Range("A6:D6").Font.Name = "Calibri"

# --- row 7: blank styled row ---------------------------------------------
$meta.Range("A7:D7").Font.Name = "Calibri"

# --- row 5: ingresos_totales (A5 left unformatted, as in the source) -----
$meta.Range("B5:D5").Font.Name = "Calibri"
$meta.Range("A5").Value = "ingresos_totales"
$meta.Range("B5").Value = " No. total de adolescentes que han ingresado al sistema SRPA en el mismo periodo y territorio. x 100"
$meta.Range("C5").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- row 6 (continued): source column, filled in after ingresos_totales --
$meta.Range("C6").Value = "Elaboración Propia"

# --- extraction date column, rows 2-6 ------------------------------------
$meta.Range("D2:D6").Value = 45722
$meta.Range("D2:D6").NumberFormat = "d-mmm-yy"

# --- rows 4 & 5 hold embedded line breaks; keep the default row height --
$meta.Rows(4).AutoFit()
$meta.Rows(5).AutoFit()

# --- make the new sheet the active / selected one, matching the author --
$meta.Range("D2:D6").Select()
$meta.Activate()
